$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md"
$ws.Range("D2").Value2 = "2016-03-23 20:55:11"
$ws.Range("A3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"
$ws.Range("D3").Value2 = "2016-03-23 20:53:40"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3c8ec4a0769d5d75b18d5368f26c3814904a30c4/e2e/a53b297b-d811-421c-9e68-f34339466385.md", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/785d70189b44bc0ae84287e30897497a8b506391/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.md") | Out-Null

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md"
$ws.Range("D2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf"
$ws.Range("E2").Value2 = "2016-03-23 20:55:06"
$ws.Range("F2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md"
$ws.Range("G2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf"
$ws.Range("H2").Value2 = "2016-03-23 20:55:33"
$ws.Range("A3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.md"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"
$ws.Range("D3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf"
$ws.Range("E3").Value2 = "2016-03-23 20:53:35"
$ws.Range("F3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.md"
$ws.Range("G3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3c8ec4a0769d5d75b18d5368f26c3814904a30c4/e2e/a53b297b-d811-421c-9e68-f34339466385.md", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3869c591d223b0be9e884b4e91550e1de6af541e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b3e6c619f64c3e23d1687d5ee8d4244991ded1a9/e2e/a53b297b-d811-421c-9e68-f34339466385.md", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d8cd7763998d15f3a5a76dcd046fb1974b0a0761/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/785d70189b44bc0ae84287e30897497a8b506391/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6fb328f726ddb727f86da016d3006d5b35d352b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b3e6c619f64c3e23d1687d5ee8d4244991ded1a9/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d8cd7763998d15f3a5a76dcd046fb1974b0a0761/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.zh-cn.xlf") | Out-Null

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md"
$ws.Range("D2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf"
$ws.Range("E2").Value2 = "2016-03-23 20:55:11"
$ws.Range("F2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md"
$ws.Range("G2").Value2 = "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf"
$ws.Range("H2").Value2 = "2016-03-23 20:55:40"
$ws.Range("A3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.md"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"
$ws.Range("D3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf"
$ws.Range("E3").Value2 = "2016-03-23 20:53:40"
$ws.Range("F3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.md"
$ws.Range("G3").Value2 = "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3c8ec4a0769d5d75b18d5368f26c3814904a30c4/e2e/a53b297b-d811-421c-9e68-f34339466385.md", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e433d3a6fc9800cfc40825cd21353444ed62014/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/374749a17b990b3e8d59587a54c6115d6b24c18a/e2e/a53b297b-d811-421c-9e68-f34339466385.md", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8553173bba9fc48d9403fc72c7610a55d957e3ef/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf", $null, $null, "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/785d70189b44bc0ae84287e30897497a8b506391/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10c75d53c48a71975083f97e1d7f68a66b707e8d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/374749a17b990b3e8d59587a54c6115d6b24c18a/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8553173bba9fc48d9403fc72c7610a55d957e3ef/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf", $null, $null, "a53b297b-d811-421c-9e68-f34339466385.4b526382c335b448856a2025d69ed8708cb473b8.de-de.xlf") | Out-Null
